$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A new daily-push row was inserted right before the old "2026/12/29" row
# (old row 605), shifting that row and everything after it down by one
# (old row 646 -> new row 647). Reproduce with a real row insert so the
# existing rows (and their values) move down intact.
$ws.Rows.Item(605).Insert()

# Fill the freshly inserted row 605 with the new day's data. Column A holds
# a "yyyy/mm/dd" label that must stay literal text (not get reinterpreted
# as a date serial by Excel's auto-detection), so format the cell as text
# first, assign the value, then restore the default/normal style so the
# saved cell carries no extra formatting - matching the plain, unstyled
# cells used throughout the rest of the column.
$ws.Cells.Item(605, 1).NumberFormat = "@"
$ws.Cells.Item(605, 1).Value = "2026/01/08"
$ws.Cells.Item(605, 1).Style = "Normal"

$ws.Cells.Item(605, 2).Value = "木"
$ws.Cells.Item(605, 3).Value = 10
$ws.Cells.Item(605, 4).Value = 200
